# Reorder the header labels in row 1, columns C:F on every worksheet so that
# "variable_trajectory_group" moves from F1 to C1, and "normalize_group",
# "trajgroup_no_vary_q", "uniform_scaling_q" each shift one column to the
# right (D1, E1, F1 respectively). Column G1 (variable_trajectory_group_
# trajectory_type) is left untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "variable_trajectory_group"
    $ws.Range("D1").Value = "normalize_group"
    $ws.Range("E1").Value = "trajgroup_no_vary_q"
    $ws.Range("F1").Value = "uniform_scaling_q"
}
